# Generate Report for Handoff
# Updates the Priority and Latest Handoff Datetime columns for the
# four "low" priority rows (rows 4-7) on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
for ($r = 4; $r -le 7; $r++) {
    $wsZh.Cells.Item($r, 5).Value = "ht"
    $wsZh.Cells.Item($r, 8).Value = "2016-10-20 10:17:36"
}

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
for ($r = 4; $r -le 7; $r++) {
    $wsDe.Cells.Item($r, 5).Value = "ht"
    $wsDe.Cells.Item($r, 8).Value = "2016-10-20 10:17:47"
}

# --- Overview sheet mirrors the de-de "Latest Handoff Datetime" value in its
# "Latest HO Xliff Generate Date" column for these same rows ---
$wsOverview = $wb.Worksheets.Item("Overview")
for ($r = 4; $r -le 7; $r++) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-10-20 10:17:47"
}
